$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Programming in R"
$ws.Range("D10").Value = "Models in R"
$ws.Range("D11").Value = "Data pipelines"

$ws.Range("D12").Select()
